$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the original formatting of an existing data row (row 2) before
# inserting, so the brand-new row can be restored to the same "plain" look
# (Insert() otherwise copies the header row's bold/centered style down).
$plainStyle = $ws.Cells.Item(2, 1).Style
$dateNumberFormat = $ws.Cells.Item(2, 4).NumberFormat

# Insert a new row above row 2, shifting existing data rows (2-6) down to (3-7)
$ws.Rows.Item(2).Insert()

# Restore the plain (non-header) style across the new row, and the date format on D
$ws.Range("A2:R2").Style = $plainStyle
$ws.Cells.Item(2, 4).NumberFormat = $dateNumberFormat

# Populate the new week's data in row 2
$ws.Cells.Item(2, 1).Value = 10
$ws.Cells.Item(2, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(2, 3).Value = "La Araucanía"
$ws.Cells.Item(2, 4).Value = 44756
$ws.Cells.Item(2, 5).Value = 9
$ws.Cells.Item(2, 6).Value = 100112036
$ws.Cells.Item(2, 7).Value = "Caigua"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 80
$ws.Cells.Item(2, 11).Value = 20000
$ws.Cells.Item(2, 12).Value = 20000
$ws.Cells.Item(2, 13).Value = 20000
$ws.Cells.Item(2, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(2, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(2, 16).Value = 1333
$ws.Cells.Item(2, 17).Value = 15
$ws.Cells.Item(2, 18).Value = "Hortaliza"
